$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet: insert a new (blank) column before column M ---
# (pushes old M:P -> N:Q, matching the "Variable Instalments" extra column
# added to the repayment schedule table)
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("M:M").Insert()

# Give the freshly inserted column a sensible custom width (close to the
# other un-"best-fit" columns on this sheet).
$wsSchedule.Columns("M:M").ColumnWidth = 7.3

# --- Sheet selection / active tab moves from "Prepay Loan" to "Repayment schedule" ---
$wsSchedule.Activate()
$wsSchedule.Range("G9").Select()
